# Appends 8 new log rows (16-23) to the "Data" sheet, matching 4 new
# measurement runs (2018.03.09 FS x2, 2018.03.28 RS x3, 2018.03.29 RS x3).
#
# Note: a handful of the new values ("2018.03.09", "2018.03.28",
# "2018.03.29") look like dates, and Excel would normally auto-convert a
# bare .Value assignment of such a string into a date serial number. To
# keep them as plain text (as the other Date-column entries in this sheet
# already are), we briefly force the cell to Text format, assign the
# string, then clear the format back off so the cell is left with the
# sheet's default (unstyled) look - only the stored value changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "2018.03.09"
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = "15:32:32"
$ws.Range("C16").Value = "FS"
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 0.1
$ws.Range("G16").Value = 0.97
$ws.Range("H16").Value = 3495
$ws.Range("I16").Value = 0.43
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 37.67676767676767
$ws.Range("L16").Value = "N/A"

$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "2018.03.09"
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = "15:52:41"
$ws.Range("C17").Value = "FS"
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 0.1
$ws.Range("G17").Value = 0.97
$ws.Range("H17").Value = 3495
$ws.Range("I17").Value = 0.21
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = 37.67676767676767
$ws.Range("L17").Value = "N/A"

$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "2018.03.28"
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = "14:40:02"
$ws.Range("C18").Value = "RS"
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 0.1
$ws.Range("G18").Value = 0.96
$ws.Range("H18").Value = 3495
$ws.Range("I18").Value = 0.15
$ws.Range("J18").Value = 1
$ws.Range("K18").Value = 32.62626262626263
$ws.Range("L18").Value = "N/A"

$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "2018.03.28"
$ws.Range("A19").ClearFormats()
$ws.Range("B19").Value = "14:40:03"
$ws.Range("C19").Value = "RS"
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 0.1
$ws.Range("G19").Value = 0.96
$ws.Range("H19").Value = 3495
$ws.Range("I19").Value = 0.14
$ws.Range("J19").Value = 1
$ws.Range("K19").Value = 32.62626262626263
$ws.Range("L19").Value = "N/A"

$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "2018.03.28"
$ws.Range("A20").ClearFormats()
$ws.Range("B20").Value = "14:40:03"
$ws.Range("C20").Value = "RS"
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 0.1
$ws.Range("G20").Value = 0.96
$ws.Range("H20").Value = 3495
$ws.Range("I20").Value = 0.14
$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 32.62626262626263
$ws.Range("L20").Value = "N/A"

$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "2018.03.29"
$ws.Range("A21").ClearFormats()
$ws.Range("B21").Value = "16:48:14"
$ws.Range("C21").Value = "RS"
$ws.Range("D21").Value = 11
$ws.Range("E21").Value = 50
$ws.Range("F21").Value = "N/A"
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 3499
$ws.Range("I21").Value = 0.15
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 27.85571142284569
$ws.Range("L21").Value = "N/A"

$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "2018.03.29"
$ws.Range("A22").ClearFormats()
$ws.Range("B22").Value = "16:48:14"
$ws.Range("C22").Value = "RS"
$ws.Range("D22").Value = 11
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = "N/A"
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 3499
$ws.Range("I22").Value = 0.13
$ws.Range("J22").Value = 1
$ws.Range("K22").Value = 27.85571142284569
$ws.Range("L22").Value = "N/A"

$ws.Range("A23").NumberFormat = "@"
$ws.Range("A23").Value = "2018.03.29"
$ws.Range("A23").ClearFormats()
$ws.Range("B23").Value = "16:48:15"
$ws.Range("C23").Value = "RS"
$ws.Range("D23").Value = 11
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = "N/A"
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 3499
$ws.Range("I23").Value = 0.13
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 27.85571142284569
$ws.Range("L23").Value = "N/A"

Write-Output "done"
